$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 154 and row 156 (columns B:G)
$tmp = $ws.Range("B154:G154").Value2
$ws.Range("B154:G154").Value2 = $ws.Range("B156:G156").Value2
$ws.Range("B156:G156").Value2 = $tmp

# Swap row 176 and row 177 (columns B:G)
$tmp = $ws.Range("B176:G176").Value2
$ws.Range("B176:G176").Value2 = $ws.Range("B177:G177").Value2
$ws.Range("B177:G177").Value2 = $tmp

# Swap row 256 and row 257 (columns B:G)
$tmp = $ws.Range("B256:G256").Value2
$ws.Range("B256:G256").Value2 = $ws.Range("B257:G257").Value2
$ws.Range("B257:G257").Value2 = $tmp

# Swap row 271 and row 272 (columns B:G)
$tmp = $ws.Range("B271:G271").Value2
$ws.Range("B271:G271").Value2 = $ws.Range("B272:G272").Value2
$ws.Range("B272:G272").Value2 = $tmp

# Swap row 338 and row 339 (columns B:G)
$tmp = $ws.Range("B338:G338").Value2
$ws.Range("B338:G338").Value2 = $ws.Range("B339:G339").Value2
$ws.Range("B339:G339").Value2 = $tmp

# Swap row 371 and row 372 (columns B:G)
$tmp = $ws.Range("B371:G371").Value2
$ws.Range("B371:G371").Value2 = $ws.Range("B372:G372").Value2
$ws.Range("B372:G372").Value2 = $tmp

# Swap row 381 and row 382 (columns B:G)
$tmp = $ws.Range("B381:G381").Value2
$ws.Range("B381:G381").Value2 = $ws.Range("B382:G382").Value2
$ws.Range("B382:G382").Value2 = $tmp

# Swap row 392 and row 393 (columns B:G)
$tmp = $ws.Range("B392:G392").Value2
$ws.Range("B392:G392").Value2 = $ws.Range("B393:G393").Value2
$ws.Range("B393:G393").Value2 = $tmp

# Swap row 423 and row 424 (columns B:G)
$tmp = $ws.Range("B423:G423").Value2
$ws.Range("B423:G423").Value2 = $ws.Range("B424:G424").Value2
$ws.Range("B424:G424").Value2 = $tmp

# Swap row 528 and row 529 (columns B:G)
$tmp = $ws.Range("B528:G528").Value2
$ws.Range("B528:G528").Value2 = $ws.Range("B529:G529").Value2
$ws.Range("B529:G529").Value2 = $tmp

# Swap row 575 and row 576 (columns B:G)
$tmp = $ws.Range("B575:G575").Value2
$ws.Range("B575:G575").Value2 = $ws.Range("B576:G576").Value2
$ws.Range("B576:G576").Value2 = $tmp

# Swap row 578 and row 579 (columns B:G)
$tmp = $ws.Range("B578:G578").Value2
$ws.Range("B578:G578").Value2 = $ws.Range("B579:G579").Value2
$ws.Range("B579:G579").Value2 = $tmp

# Swap row 582 and row 583 (columns B:G)
$tmp = $ws.Range("B582:G582").Value2
$ws.Range("B582:G582").Value2 = $ws.Range("B583:G583").Value2
$ws.Range("B583:G583").Value2 = $tmp

# Swap row 585 and row 586 (columns B:G)
$tmp = $ws.Range("B585:G585").Value2
$ws.Range("B585:G585").Value2 = $ws.Range("B586:G586").Value2
$ws.Range("B586:G586").Value2 = $tmp

# Swap row 591 and row 592 (columns B:G)
$tmp = $ws.Range("B591:G591").Value2
$ws.Range("B591:G591").Value2 = $ws.Range("B592:G592").Value2
$ws.Range("B592:G592").Value2 = $tmp

# Swap row 679 and row 680 (columns B:G)
$tmp = $ws.Range("B679:G679").Value2
$ws.Range("B679:G679").Value2 = $ws.Range("B680:G680").Value2
$ws.Range("B680:G680").Value2 = $tmp

# Swap row 701 and row 702 (columns B:G)
$tmp = $ws.Range("B701:G701").Value2
$ws.Range("B701:G701").Value2 = $ws.Range("B702:G702").Value2
$ws.Range("B702:G702").Value2 = $tmp

# Swap row 864 and row 865 (columns B:G)
$tmp = $ws.Range("B864:G864").Value2
$ws.Range("B864:G864").Value2 = $ws.Range("B865:G865").Value2
$ws.Range("B865:G865").Value2 = $tmp
